$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Select the original range that held the data and clear its contents,
# leaving cell formatting (styles) untouched.
$ws.Range("A1:E11").Select()
$ws.Range("A1:E11").ClearContents()

$ws.Range("A1:E5").EntireRow.AutoFit()

# Final selection ends up on E4 (single cell) as in the target workbook.
$ws.Range("E4").Select()
